$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels in A1 and B1
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"

# Update the one-hot block order values (rows 2-7, columns A-F)
$values = @(
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,0,1,0,0),
    @(0,0,0,0,1,0),
    @(0,0,1,0,0,0),
    @(1,0,0,0,0,0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $rowIndex = $i + 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowVals[$j]
    }
}
